$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58 (pushes the former rows 58-74 down to 59-75,
# carrying their data/formatting with them, matching Excel's native Insert
# behaviour used when a fresh weekly record is prepended to the table).
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new weekly price record.
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = 'Vega Modelo de Temuco'
$ws.Range("C58").Value = 'La Araucanía'
$ws.Range("D58").Value = 44642
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = 'Fruta'
$ws.Range("G58").Value = 100108
$ws.Range("H58").Value = 'Tropicales y subtropicales'
$ws.Range("I58").Value = 100108004
$ws.Range("J58").Value = 'Papaya'
$ws.Range("K58").Value = 'Cultivar IV Región'
$ws.Range("L58").Value = 'Primera'
$ws.Range("M58").Value = 30
$ws.Range("N58").Value = 24000
$ws.Range("O58").Value = 24000
$ws.Range("P58").Value = 24000
$ws.Range("Q58").Value = '$/bandeja 10 kilos'
$ws.Range("R58").Value = 'Provincia del Elquí'
$ws.Range("S58").Value = 2400
$ws.Range("T58").Value = 10
